$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = "42.736.57"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "2.279.56"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "'310.65"
$ws.Range("E5").Value = "  -2.75%  "
$ws.Range("D6").Value = "'103.74"
$ws.Range("E6").Value = "  +3.10%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "'1.01"
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").Value = "'38.68"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").Value = "'0.0898"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "'8.21"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").Value = "'0.970"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "'14.98"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "2.628.12"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "2.286.84"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "42.367.13"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").Value = "'7.23"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").Value = "'13.33"
$ws.Range("E21").Value = "  +5.24%  "
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("E23").Value = "  -3.44%  "
$ws.Range("D24").Value = "'262.58"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").Value = "'10.63"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D30").Value = "'22.14"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").Value = "'35.57"
$ws.Range("E31").Value = "  -4.74%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").Value = "'0.0853"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("D37").Value = "'4.47"
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D38").Value = "'0.0347"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").Value = "'3.71"
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("D40").Value = "'2.70"
$ws.Range("E40").Value = "  -2.30%  "
$ws.Range("D41").Value = "'1.55"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").Value = "'97.52"
$ws.Range("E42").Value = "  +5.96%  "
$ws.Range("D43").Value = "'68.51"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D46").Value = "1.719.47"
$ws.Range("E46").Value = "  +6.96%  "
$ws.Range("D47").Value = "'11.84"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "'109.62"
$ws.Range("E48").Value = "  -4.70%  "
$ws.Range("D49").Value = "'76.06"
$ws.Range("E49").Value = "  -3.44%  "
$ws.Range("D50").Value = "'5.17"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").Value = "'8.59"
$ws.Range("E51").Value = "  -3.67%  "

# --- Row swaps: Filecoin <-> Toncoin (rows 28/29) ---
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.32"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'6.92"
$ws.Range("E29").Value = "  +14.68%  "

# --- Row swaps: Algorand <-> FirstDigitalUSD (rows 44/45) ---
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.225"
$ws.Range("E45").Value = "  +0.78%  "
